$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "Warrant Cancelled Report" sheet by copying the
#    existing "Warrant Accepted Report" sheet (same mapping layout),
#    placed right before "Warrant Rejected Report".
# ------------------------------------------------------------------
$wsAccepted = $wb.Worksheets.Item("Warrant Accepted Report")
$wsRejected = $wb.Worksheets.Item("Warrant Rejected Report")

$wsAccepted.Copy($wsRejected)

$wsCancelled = $wb.Worksheets.Item("Warrant Accepted Report (2)")
$wsCancelled.Name = "Warrant Cancelled Report"

# ------------------------------------------------------------------
# 2. Update the title cell.
# ------------------------------------------------------------------
$wsCancelled.Range("B1").Value = "Warrant Cancelled Report"

# ------------------------------------------------------------------
# 3. Remove the two rows that don't apply to this report:
#      row 13 - "Law Enforcement ORI"
#      row 16 - "Warrant Entry Category Code Text"
#    Delete the higher-numbered row first so the other row index
#    stays valid.
# ------------------------------------------------------------------
$wsCancelled.Rows("16").Delete()
$wsCancelled.Rows("13").Delete()

# ------------------------------------------------------------------
# 4. Update the XPath-mapping text in column F for every remaining
#    data row so it refers to "wcr-doc:WarrantAcceptedReport" (the
#    Warrant Cancelled Report schema) instead of the original
#    "war-doc:WarrantAcceptedReport" text.
# ------------------------------------------------------------------
$wsCancelled.Range("F4").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonGivenName"
$wsCancelled.Range("F5").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonMiddleName"
$wsCancelled.Range("F6").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonSurName"
$wsCancelled.Range("F7").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonNameSuffixText"
$wsCancelled.Range("F8").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonBirthDate/nc:Date"
$wsCancelled.Range("F9").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonRaceText"
$wsCancelled.Range("F10").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonSexText"
$wsCancelled.Range("F13").Value = "wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderEnforcementAgency/wir-ext:AgencyRecordIdentification/nc:IdentificationID"
$wsCancelled.Range("F14").Value = "wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderRequestEntity/nc:EntityPerson/wir-ext:PersonEmployeeIdentification/nc:IdentificationID"
$wsCancelled.Range("F15").Value = "wcr-doc:WarrantAcceptedReport/j:Warrant/wir-ext:WarrantAugmentation/wir-ext:StateWarrantRepositoryIdentification/nc:IdentificationID/#text"

# ------------------------------------------------------------------
# 5. View tweaks to match the authored workbook: the new sheet is the
#    active/selected tab, scrolled/selected at A2.
# ------------------------------------------------------------------
$wsCancelled.Range("A2").Select()
$wsCancelled.Activate()

$wb.Windows.Item(1).ActiveSheet.Application.ActiveWindow.ScrollRow = 7
